$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = '28.034.94'
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  +2.04%  '
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = '1.648.12'
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  +1.88%  '
$ws.Cells.Item(4, 5).Value = '  +0.03%  '
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '213.96'
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +1.30%  '
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '0.528'
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +0.52%  '
$ws.Cells.Item(7, 5).Value = '  +0.01%  '
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = '23.58'
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  +3.37%  '
$ws.Cells.Item(9, 5).Value = '  +1.70%  '
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = '0.0616'
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +0.59%  '
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = '0.0873'
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -1.57%  '
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '1.883.28'
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +1.94%  '
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = '1.649.83'
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +1.91%  '
$ws.Cells.Item(14, 5).Value = '  +1.22%  '
$ws.Cells.Item(15, 5).Value = '  +3.20%  '
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = '65.72'
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +1.24%  '
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = '28.049.34'
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +2.15%  '
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = '232.60'
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +0.02%  '
$ws.Cells.Item(19, 2).Value = 'Chainlink'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '7.68'
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +1.73%  '
$ws.Cells.Item(20, 2).Value = 'ShibaInu'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '0.0₃0724'
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +0.94%  '
$ws.Cells.Item(21, 5).Value = '  +0.00%  '
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '10.73'
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +6.04%  '
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '4.40'
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +2.44%  '
$ws.Cells.Item(24, 5).Value = '  +3.51%  '
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = '152.30'
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +1.21%  '
$ws.Cells.Item(26, 5).Value = '  +1.11%  '
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = '15.79'
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +1.61%  '
$ws.Cells.Item(29, 5).Value = '  -0.07%  '
$ws.Cells.Item(30, 5).Value = '  +1.34%  '
$ws.Cells.Item(31, 5).Value = '  +0.32%  '
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '3.34'
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +2.20%  '
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = '1.446.02'
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -1.72%  '
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = '3.08'
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +0.30%  '
$ws.Cells.Item(35, 5).Value = '  +1.52%  '
$ws.Cells.Item(36, 5).Value = '  -0.21%  '
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = '0.892'
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +3.68%  '
$ws.Cells.Item(38, 5).Value = '  +1.54%  '
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '0.933'
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -2.18%  '
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '0.558'
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -0.01%  '
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '69.40'
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +1.98%  '
$ws.Cells.Item(42, 5).Value = '  +3.62%  '
$ws.Cells.Item(44, 5).Value = '  -0.03%  '
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '1.84'
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +5.85%  '
$ws.Cells.Item(46, 5).Value = '  +2.89%  '
$ws.Cells.Item(47, 5).Value = '  +0.94%  '
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '1.791.74'
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +1.82%  '
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = '89.22'
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +2.97%  '
$ws.Cells.Item(50, 5).Value = '  -0.46%  '
$ws.Cells.Item(51, 5).Value = '  +0.38%  '
